# Apply crypto price/volume updates per the commit diff
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "24.473.38"
$ws.Range("E2").Value = "  +9.03%  "

$ws.Range("D3").Value = "1.682.07"
$ws.Range("E3").Value = "  +4.40%  "

$ws.Range("E4").Value = "  +0.11%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "307.68"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.99%  "

$ws.Range("E6").Value = "  +0.68%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.3710"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +0.40%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3448"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +0.86%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "47.88"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +13.06%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "1.177"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +3.34%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07277"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +2.61%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "1.001"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +0.03%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "6.134"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +3.27%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "20.35"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +2.70%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "6.760"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +1.34%  "

$ws.Range("D16").Value = "1.681.16"
$ws.Range("E16").Value = "  +4.87%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.00001109"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +1.97%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.9985"
$ws.Range("D18").Style = "Normal"

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.06685"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -1.66%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "81.28"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +3.62%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "16.49"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +2.28%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "6.129"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +1.12%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "12.21"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +2.69%  "

$ws.Range("D24").Value = "24.432.83"
$ws.Range("E24").Value = "  +8.74%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.449"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +2.14%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.671"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +4.57%  "

$ws.Range("B27").Value = "Monero"
$ws.Range("C27").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "153.84"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +2.17%  "

$ws.Range("B28").Value = "EthereumClassic"
$ws.Range("C28").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "19.51"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -0.55%  "

$ws.Range("B29").Value = "WrappedliquidstakedEther2.0"
$ws.Range("C29").Value = "https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth"
$ws.Range("D29").Value = "1.867.83"
$ws.Range("E29").Value = "  +4.66%  "

$ws.Range("B30").Value = "BitcoinCash"
$ws.Range("C30").Value = "https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "127.37"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +3.78%  "

$ws.Range("B31").Value = "Filecoin"
$ws.Range("C31").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "6.342"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +3.27%  "

$ws.Range("B32").Value = "HuobiToken"
$ws.Range("C32").Value = "https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "4.059"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +1.40%  "

$ws.Range("B33").Value = "ImmutableX"
$ws.Range("C33").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.9799"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +2.65%  "

$ws.Range("B34").Value = "Stellar"
$ws.Range("C34").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.08464"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +2.16%  "

$ws.Range("B35").Value = "WEMIXTOKEN"
$ws.Range("C35").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.708"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +3.71%  "

$ws.Range("B36").Value = "Aptos"
$ws.Range("C36").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "12.43"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +3.25%  "

$ws.Range("B37").Value = "Hedera"
$ws.Range("C37").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.06504"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +6.52%  "

$ws.Range("B38").Value = "InternetComputer(DFINITY)"
$ws.Range("C38").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "5.379"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +2.05%  "

$ws.Range("B39").Value = "FraxShare"
$ws.Range("C39").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "8.877"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +2.85%  "

$ws.Range("B40").Value = "VeChain"
$ws.Range("C40").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.02331"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +4.54%  "

$ws.Range("B41").Value = "TrustWalletToken"
$ws.Range("C41").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "1.257"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -0.91%  "

$ws.Range("B42").Value = "Algorand"
$ws.Range("C42").Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.2117"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +4.19%  "

$ws.Range("B43").Value = "TheSandbox"
$ws.Range("C43").Value = "https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.6165"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +3.75%  "

$ws.Range("B44").Value = "Frax"
$ws.Range("C44").Value = "https://coinranking.com/coin/KfWtaeV1W+frax-frax"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.9980"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +0.62%  "

$ws.Range("B45").Value = "PancakeSwap"
$ws.Range("C45").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "3.773"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -1.38%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "13.11"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -0.49%  "

$ws.Range("B47").Value = "Decentraland"
$ws.Range("C47").Value = "https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.5952"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +3.98%  "

$ws.Range("B48").Value = "Quant"
$ws.Range("C48").Value = "https://coinranking.com/coin/bauj_21eYVwso+quant-qnt"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "127.46"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -0.02%  "

$ws.Range("B49").Value = "NEARProtocol"
$ws.Range("C49").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "2.027"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +2.10%  "

$ws.Range("B50").Value = "Cronos"
$ws.Range("C50").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.07203"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +5.68%  "

$ws.Range("B51").Value = "Aave"
$ws.Range("C51").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "76.13"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +2.65%  "
